$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 holds three "blog" card entries (columns C, E, I) whose text ends
# with "ser: <n>". A new blog post (ser: 134) was added, bumping each
# existing reference up by one: 131 -> 132 (I8), 132 -> 133 (E8),
# 133 -> 134 (C8).
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 134"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 133"
$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 132"
